$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 19 (pushes old rows 19+ down, including the signature block at 23-24 -> 24-25)
$ws.Rows("19:19").Insert()

# The newly inserted row 19 inherits row-18's old "last row" formatting; copy it down
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)  # xlPasteFormats

# Row 18 is no longer the last data row -> give it the "middle" row formatting (same as rows 16/17)
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)  # xlPasteFormats

# Re-order the "Periodo Mora" column to ascending (2505, 2506, 2507) and
# populate the new data row with the next period (2508)
$ws.Range("E16").Value = "2505"
$ws.Range("E17").Value = "2506"
$ws.Range("E18").Value = "2507"

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1041945903"
$ws.Range("D19").Value = "LUIS FERNANDO GUARDO CORREA"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 64000
$ws.Range("G19").Value = 1600000

# Update Valor Mora total
$ws.Range("E11").Value = 256000

# Update Cant. Periodos
$ws.Range("F13").Value = 4

$excel.CutCopyMode = $false
